$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.989.23'
$ws.Range("E2").Value = '  +2.49%  '

$ws.Range("D3").Value = '2.254.21'
$ws.Range("E3").Value = '  +1.70%  '

$ws.Range("E4").Value = '  +0.16%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '319.44'
$ws.Range("E5").Value = '  +0.07%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '101.68'
$ws.Range("E6").Value = '  +3.33%  '

$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("E8").Value = '  +0.13%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.553'
$ws.Range("E9").Value = '  -0.36%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '37.34'
$ws.Range("E10").Value = '  +2.07%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0831'
$ws.Range("E11").Value = '  +0.99%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.60'
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("E13").Value = '  -1.27%  '

$ws.Range("D14").Value = '2.599.04'
$ws.Range("E14").Value = '  +1.73%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.857'
$ws.Range("E15").Value = '  -0.30%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '14.44'
$ws.Range("E16").Value = '  +1.04%  '

$ws.Range("D17").Value = '2.255.40'
$ws.Range("E17").Value = '  +1.81%  '

$ws.Range("D18").Value = '43.884.03'
$ws.Range("E18").Value = '  +2.55%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.32'
$ws.Range("E19").Value = '  -4.21%  '

$ws.Range("D20").Value = '0.0₃0985'
$ws.Range("E20").Value = '  +2.54%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.52'
$ws.Range("E21").Value = '  +0.14%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '65.76'
$ws.Range("E22").Value = '  +1.48%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '3.15'
$ws.Range("E23").Value = '  -0.75%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '235.59'
$ws.Range("E24").Value = '  +0.03%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.10'
$ws.Range("E25").Value = '  -1.96%  '

$ws.Range("E26").Value = '  -0.10%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.20'
$ws.Range("E27").Value = '  +2.82%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.19'
$ws.Range("E28").Value = '  -2.15%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '37.83'
$ws.Range("E29").Value = '  +6.54%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.19'
$ws.Range("E30").Value = '  -1.78%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '161.29'
$ws.Range("E31").Value = '  +6.28%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '20.22'
$ws.Range("E32").Value = '  +0.23%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0851'
$ws.Range("E33").Value = '  -0.94%  '

$ws.Range("E34").Value = '  +1.72%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.115'
$ws.Range("E35").Value = '  +11.76%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.93'
$ws.Range("E36").Value = '  +1.86%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.02'
$ws.Range("E37").Value = '  -4.39%  '

$ws.Range("E38").Value = '  -1.44%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '16.75'
$ws.Range("E39").Value = '  +23.88%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.72'
$ws.Range("E40").Value = '  +2.11%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.21'
$ws.Range("E41").Value = '  -4.03%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0315'
$ws.Range("E42").Value = '  -1.18%  '

$ws.Range("E43").Value = '  +0.20%  '

$ws.Range("D44").Value = '1.794.54'
$ws.Range("E44").Value = '  +3.90%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.199'
$ws.Range("E45").Value = '  -1.74%  '

$ws.Range("B46").Value = 'ordi'
$ws.Range("C46").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '75.51'
$ws.Range("E46").Value = '  +2.30%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '82.59'
$ws.Range("E47").Value = '  -1.72%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '5.21'

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '104.84'
$ws.Range("E49").Value = '  +2.28%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.69'
$ws.Range("E50").Value = '  +9.63%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '58.15'
$ws.Range("E51").Value = '  +1.45%  '
